$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scen_2")

$ws.Range("B1").Comment.Delete()
$c = $ws.Range("B3").AddComment("Tom Tidhar:`nProportion of diagnosed PLHIV who are on treatment")
Write-Host "Author before:" $c.Author
$c.Author = "Tom Tidhar"
Write-Host "Author after:" $c.Author
Write-Host "Text:" $c.Text()
